$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 443.1875
$ws.Range("I28").Value = 406.06668
$ws.Range("J28").Value = 1000
$ws.Range("K28").Value = 406.06668
$ws.Range("L28").Value = 1000
$ws.Range("M28").Value = 78.93331999999998
$ws.Range("N28").Value = -1970
$ws.Range("H40").Value = 2020
$ws.Range("I40").Value = 742.8570999999999
$ws.Range("K40").Value = 742.8570999999999
$ws.Range("M40").Value = -567.8570999999999
$ws.Range("H62").Value = 4002.6
$ws.Range("I62").Value = 3036.8333
$ws.Range("K62").Value = 3036.8333
$ws.Range("M62").Value = -2412.8333
$ws.Range("H65").Value = 4002.6
$ws.Range("I65").Value = 3036.8333
$ws.Range("K65").Value = 15184.1665
$ws.Range("M65").Value = -12064.1665
$ws.Range("H113").Value = 4336.1665
$ws.Range("I113").Value = 4099.875
$ws.Range("J113").Value = 4525.2
$ws.Range("K113").Value = 4099.875
$ws.Range("L113").Value = 4525.2
$ws.Range("M113").Value = -845.875
$ws.Range("N113").Value = -11033.2
$ws.Range("H116").Value = 4256.724
$ws.Range("I116").Value = 3551.7856
$ws.Range("K116").Value = 3551.7856
$ws.Range("M116").Value = -109.7856000000002
$ws.Range("H137").Value = 3852039.8
$ws.Range("I137").Value = 5889210
$ws.Range("J137").Value = 4051.7778
$ws.Range("K137").Value = 17667630
$ws.Range("L137").Value = 12155.3334
$ws.Range("M137").Value = -17665080
$ws.Range("N137").Value = -17255.3334
$ws.Range("H138").Value = 4700.471
$ws.Range("I138").Value = 2861.4614
$ws.Range("J138").Value = 5787.159
$ws.Range("K138").Value = 8584.3842
$ws.Range("L138").Value = 17361.477
$ws.Range("M138").Value = -3444.3842
$ws.Range("N138").Value = -27641.477

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 17938
$ws.Range("J37").Value = 17938
$ws.Range("L37").Value = 17938
$ws.Range("N37").Value = -18484
$ws.Range("H44").Value = 21308.166
$ws.Range("J44").Value = 21308.166
$ws.Range("L44").Value = 21308.166
$ws.Range("N44").Value = -22284.166
$ws.Range("H45").Value = 1335.225
$ws.Range("I45").Value = 1063.0264
$ws.Range("K45").Value = 1063.0264
$ws.Range("M45").Value = -686.0264
$ws.Range("H55").Value = 24963.25
$ws.Range("J55").Value = 24963.25
$ws.Range("L55").Value = 24963.25
$ws.Range("N55").Value = -25593.25
$ws.Range("H61").Value = 2432.8086
$ws.Range("I61").Value = 1804.1765
$ws.Range("J61").Value = 4076.923
$ws.Range("K61").Value = 1804.1765
$ws.Range("L61").Value = 4076.923
$ws.Range("M61").Value = -1592.1765
$ws.Range("N61").Value = -4500.923
$ws.Range("H80").Value = 29591.572
$ws.Range("J80").Value = 29591.572
$ws.Range("L80").Value = 29591.572
$ws.Range("N80").Value = -31587.572
$ws.Range("H83").Value = 29591.572
$ws.Range("J83").Value = 29591.572
$ws.Range("L83").Value = 88774.716
$ws.Range("N83").Value = -98758.716
$ws.Range("H97").Value = 447.64706
$ws.Range("I97").Value = 466.25
$ws.Range("K97").Value = 466.25
$ws.Range("M97").Value = 29.75
$ws.Range("H132").Value = 2911.5
$ws.Range("I132").Value = 2341.647
$ws.Range("J132").Value = 3656.6924
$ws.Range("K132").Value = 7024.941
$ws.Range("L132").Value = 10970.0772
$ws.Range("M132").Value = -4494.941
$ws.Range("N132").Value = -16030.0772
$ws.Range("H136").Value = 2432.8086
$ws.Range("I136").Value = 1804.1765
$ws.Range("J136").Value = 4076.923
$ws.Range("K136").Value = 5412.529500000001
$ws.Range("L136").Value = 12230.769
$ws.Range("M136").Value = -2862.529500000001
$ws.Range("N136").Value = -17330.769

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 729.4074000000001
$ws.Range("I94").Value = 570.3043
$ws.Range("J94").Value = 1644.25
$ws.Range("K94").Value = 570.3043
$ws.Range("L94").Value = 1644.25
$ws.Range("M94").Value = -119.3043
$ws.Range("N94").Value = -2546.25
$ws.Range("H105").Value = 2430.1667
$ws.Range("I105").Value = 2177.5
$ws.Range("J105").Value = 2935.5
$ws.Range("K105").Value = 2177.5
$ws.Range("L105").Value = 2935.5
$ws.Range("M105").Value = -430.5
$ws.Range("N105").Value = -6429.5
$ws.Range("H134").Value = 2897.0378
$ws.Range("I134").Value = 2609.152
$ws.Range("J134").Value = 4788.857
$ws.Range("K134").Value = 7827.456
$ws.Range("L134").Value = 14366.571
$ws.Range("M134").Value = -5292.456
$ws.Range("N134").Value = -19436.571

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1009.94116
$ws.Range("I22").Value = 382.5
$ws.Range("K22").Value = 382.5
$ws.Range("M22").Value = -32.5
$ws.Range("H31").Value = 1855389.8
$ws.Range("I31").Value = 2779827.2
$ws.Range("J31").Value = 6514.778
$ws.Range("K31").Value = 2779827.2
$ws.Range("L31").Value = 6514.778
$ws.Range("M31").Value = -2779532.2
$ws.Range("N31").Value = -7104.778
$ws.Range("H34").Value = 1855389.8
$ws.Range("I34").Value = 2779827.2
$ws.Range("J34").Value = 6514.778
$ws.Range("K34").Value = 2779827.2
$ws.Range("L34").Value = 6514.778
$ws.Range("M34").Value = -2779625.2
$ws.Range("N34").Value = -6918.778
$ws.Range("H58").Value = 9807171
$ws.Range("I58").Value = 1768.4445
$ws.Range("J58").Value = 20838250
$ws.Range("K58").Value = 1768.4445
$ws.Range("L58").Value = 20838250
$ws.Range("M58").Value = -1565.4445
$ws.Range("N58").Value = -20838656
$ws.Range("H59").Value = 21438.428
$ws.Range("J59").Value = 21438.428
$ws.Range("L59").Value = 21438.428
$ws.Range("N59").Value = -23728.428
$ws.Range("H60").Value = 14367.667
$ws.Range("H74").Value = 24309.572
$ws.Range("J74").Value = 24309.572
$ws.Range("L74").Value = 24309.572
$ws.Range("N74").Value = -26057.572
$ws.Range("H77").Value = 24309.572
$ws.Range("J77").Value = 24309.572
$ws.Range("L77").Value = 72928.716
$ws.Range("N77").Value = -81664.716
$ws.Range("H107").Value = 1279.5143
$ws.Range("I107").Value = 1154.8
$ws.Range("J107").Value = 1591.3
$ws.Range("K107").Value = 1154.8
$ws.Range("L107").Value = 1591.3
$ws.Range("M107").Value = 765.2
$ws.Range("N107").Value = -5431.3
$ws.Range("H133").Value = 21289.5
$ws.Range("J133").Value = 21289.5
$ws.Range("L133").Value = 21289.5
$ws.Range("N133").Value = -26349.5
$ws.Range("H135").Value = 38590
$ws.Range("J135").Value = 38590
$ws.Range("L135").Value = 38590
$ws.Range("N135").Value = -48730
$ws.Range("H136").Value = 9807171
$ws.Range("I136").Value = 1768.4445
$ws.Range("J136").Value = 20838250
$ws.Range("K136").Value = 5305.333500000001
$ws.Range("L136").Value = 62514750
$ws.Range("M136").Value = -2755.333500000001
$ws.Range("N136").Value = -62519850
$ws.Range("H138").Value = 27104.16
$ws.Range("J138").Value = 27104.16
$ws.Range("L138").Value = 27104.16
$ws.Range("N138").Value = -37384.16
$ws.Range("H140").Value = 60000
$ws.Range("J140").Value = 60000
$ws.Range("L140").Value = 60000
$ws.Range("N140").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 8533.666999999999
$ws.Range("I87").Value = 4813
$ws.Range("K87").Value = 14439
$ws.Range("M87").Value = -13191
$ws.Range("H90").Value = 8533.666999999999
$ws.Range("I90").Value = 4813
$ws.Range("K90").Value = 43317
$ws.Range("M90").Value = -37077
$ws.Range("H131").Value = 1455.2593
$ws.Range("J131").Value = 1222.7675
$ws.Range("L131").Value = 3668.3025
$ws.Range("N131").Value = -13748.3025

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 87.304344
$ws.Range("I2").Value = 102.10526
$ws.Range("J2").Value = 17
$ws.Range("K2").Value = 102.10526
$ws.Range("L2").Value = 17
$ws.Range("M2").Value = 10.89474
$ws.Range("N2").Value = -243
$ws.Range("H103").Value = 23960.8
$ws.Range("I103").Value = 14750
$ws.Range("J103").Value = 30101.334
$ws.Range("K103").Value = 14750
$ws.Range("L103").Value = 30101.334
$ws.Range("M103").Value = -13578
$ws.Range("N103").Value = -32445.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2717.4
$ws.Range("I93").Value = 1982.1428
$ws.Range("K93").Value = 1982.1428
$ws.Range("M93").Value = -734.1428000000001
$ws.Range("H122").Value = 5860
$ws.Range("I122").Value = 3900
$ws.Range("J122").Value = 7166.6665
$ws.Range("K122").Value = 11700
$ws.Range("L122").Value = 21499.9995
$ws.Range("M122").Value = -9250
$ws.Range("N122").Value = -26399.9995
$ws.Range("H132").Value = 2508.7844
$ws.Range("I132").Value = 1742.125
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 5226.375
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -2696.375
$ws.Range("N132").Value = -16460

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10759
$ws.Range("J41").Value = 10759
$ws.Range("L41").Value = 10759
$ws.Range("N41").Value = -11539
$ws.Range("H81").Value = 1148.4166
$ws.Range("I81").Value = 1148.4166
$ws.Range("K81").Value = 2296.8332
$ws.Range("M81").Value = -1235.8332
$ws.Range("H84").Value = 1148.4166
$ws.Range("I84").Value = 1148.4166
$ws.Range("K84").Value = 11484.166
$ws.Range("M84").Value = -6180.166000000001
$ws.Range("H132").Value = 1897574.2
$ws.Range("I132").Value = 2225251.8
$ws.Range("K132").Value = 6675755.399999999
$ws.Range("M132").Value = -6673225.399999999
$ws.Range("H135").Value = 80143
$ws.Range("J135").Value = 80143
$ws.Range("L135").Value = 80143
$ws.Range("N135").Value = -90283
$ws.Range("H136").Value = 2346.3865
$ws.Range("I136").Value = 1727
$ws.Range("K136").Value = 5181
$ws.Range("M136").Value = -2631
$ws.Range("H138").Value = 29429
$ws.Range("J138").Value = 29429
$ws.Range("L138").Value = 29429
$ws.Range("N138").Value = -39709
